$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Set the Approved/Rejected value for row 8 (TestScenario_2.TestCase_1)
$ws.Range("I8").Value = "Approved"

# Update the active selection to match the saved view state
$ws.Range("E8").Select()
